# Generate Report for Handoff
# Update the "591e6bd7-243e-4d56-b3dd-c8d740216bf9" entry (row 7 on each sheet)
# with fresh handoff/handback timestamps.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: column D = "Latest Handoff Date", row 7 is the 591e6bd7... entry
$overview.Range("D7").Value = "2016-28-17 12:28:52"

# zh-cn sheet: column E = "Latest Handoff Datetime", row 7 is the 591e6bd7... entry
$zhcn.Range("E7").Value = "2016-03-17 12:28:48"

# de-de sheet: column E = "Latest Handoff Datetime" (used here for Handback timestamp), row 7 is the 591e6bd7... entry
$dede.Range("E7").Value = "2016-03-17 12:28:52"
